{"js": "const pairs = [\n  [\"63\u00f77=9, 0\", \"44\u00f77=6, 2\"],\n  [\"63\u00f79=7, 0\", \"95\u00f78=11, 7\"],\n  [\"31\u00f77=4, 3\", \"24\u00f79=2, 6\"],\n  [\"20\u00f75=4, 0\", \"25\u00f72=12, 1\"],\n  [\"82\u00f76=13, 4\", \"61\u00f76=10, 1\"],\n  [\"11\u00f79=1, 2\", \"73\u00f75=14, 3\"],\n  [\"24\u00f72=12, 0\", \"92\u00f77=13, 1\"],\n  [\"36\u00f78=4, 4\", \"76\u00f74=19, 0\"],\n  [\"78\u00f75=15, 3\", \"73\u00f76=12, 1\"],\n  [\"64\u00f73=21, 1\", \"96\u00f79=10, 6\"],\n  [\"78\u00f78=9, 6\", \"20\u00f73=6, 2\"],\n  [\"19\u00f76=3, 1\", \"89\u00f74=22, 1\"],\n  [\"65\u00f79=7, 2\", \"70\u00f77=10, 0\"],\n  [\"20\u00f77=2, 6\", \"69\u00f79=7, 6\"],\n  [\"17\u00f77=2, 3\", \"56\u00f75=11, 1\"],\n  [\"17\u00f73=5, 2\", \"45\u00f75=9, 0\"],\n  [\"64\u00f72=32, 0\", \"32\u00f77=4, 4\"],\n  [\"33\u00f75=6, 3\", \"46\u00f72=23, 0\"],\n  [\"61\u00f72=30, 1\", \"34\u00f76=5, 4\"],\n  [\"50\u00f73=16, 2\", \"44\u00f74=11, 0\"],\n  [\"39\u00f75=7, 4\", \"68\u00f78=8, 4\"],\n  [\"32\u00f74=8, 0\", \"75\u00f72=37, 1\"],\n  [\"78\u00f76=13, 0\", \"81\u00f73=27, 0\"],\n  [\"47\u00f73=15, 2\", \"64\u00f77=9, 1\"],\n  [\"77\u00f75=15, 2\", \"72\u00f73=24, 0\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('text');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"63\u00f77=9, 0\", \"44\u00f77=6, 2\")\n    ,@(\"63\u00f79=7, 0\", \"95\u00f78=11, 7\")\n    ,@(\"31\u00f77=4, 3\", \"24\u00f79=2, 6\")\n    ,@(\"20\u00f75=4, 0\", \"25\u00f72=12, 1\")\n    ,@(\"82\u00f76=13, 4\", \"61\u00f76=10, 1\")\n    ,@(\"11\u00f79=1, 2\", \"73\u00f75=14, 3\")\n    ,@(\"24\u00f72=12, 0\", \"92\u00f77=13, 1\")\n    ,@(\"36\u00f78=4, 4\", \"76\u00f74=19, 0\")\n    ,@(\"78\u00f75=15, 3\", \"73\u00f76=12, 1\")\n    ,@(\"64\u00f73=21, 1\", \"96\u00f79=10, 6\")\n    ,@(\"78\u00f78=9, 6\", \"20\u00f73=6, 2\")\n    ,@(\"19\u00f76=3, 1\", \"89\u00f74=22, 1\")\n    ,@(\"65\u00f79=7, 2\", \"70\u00f77=10, 0\")\n    ,@(\"20\u00f77=2, 6\", \"69\u00f79=7, 6\")\n    ,@(\"17\u00f77=2, 3\", \"56\u00f75=11, 1\")\n    ,@(\"17\u00f73=5, 2\", \"45\u00f75=9, 0\")\n    ,@(\"64\u00f72=32, 0\", \"32\u00f77=4, 4\")\n    ,@(\"33\u00f75=6, 3\", \"46\u00f72=23, 0\")\n    ,@(\"61\u00f72=30, 1\", \"34\u00f76=5, 4\")\n    ,@(\"50\u00f73=16, 2\", \"44\u00f74=11, 0\")\n    ,@(\"39\u00f75=7, 4\", \"68\u00f78=8, 4\")\n    ,@(\"32\u00f74=8, 0\", \"75\u00f72=37, 1\")\n    ,@(\"78\u00f76=13, 0\", \"81\u00f73=27, 0\")\n    ,@(\"47\u00f73=15, 2\", \"64\u00f77=9, 1\")\n    ,@(\"77\u00f75=15, 2\", \"72\u00f73=24, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
